$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) updates - apostrophe-prefixed to preserve text type / exact formatting
$ws.Range("D2").Value = "'275.92"
$ws.Range("D3").Value = "'21.06"
$ws.Range("D5").Value = "'0.06184"
$ws.Range("D6").Value = "'3.578"
$ws.Range("D7").Value = "'1.524"
$ws.Range("D9").Value = "'0.8226"
$ws.Range("D10").Value = "'0.1640"
$ws.Range("D11").Value = "'0.08219"
$ws.Range("D12").Value = "'0.03435"
$ws.Range("D13").Value = "'0.03128"
$ws.Range("D14").Value = "'0.09125"
$ws.Range("D15").Value = "'3.770"
$ws.Range("D16").Value = "'0.001628"
$ws.Range("D17").Value = "'0.04698"
$ws.Range("D18").Value = "'0.006436"
$ws.Range("D19").Value = "'0.006141"
$ws.Range("D20").Value = "'0.001067"
$ws.Range("D22").Value = "'3.727"
$ws.Range("D23").Value = "'2.315"
$ws.Range("D24").Value = "'0.01385"
$ws.Range("D26").Value = "'0.1231"
$ws.Range("D28").Value = "'0.0002736"
$ws.Range("D41").Value = "'0.005401"
$ws.Range("D42").Value = "'0.007057"
$ws.Range("D43").Value = "'0.1103"
$ws.Range("D44").Value = "'0.01116"
$ws.Range("D45").Value = "'0.00006264"
$ws.Range("D47").Value = "'0.8450"
$ws.Range("D48").Value = "'0.001385"
$ws.Range("D49").Value = "'0.00001900"

# Column E (Volume/1h label) updates
$ws.Range("E19").Value = "18HotbitTokenHTBBestin24h"
$ws.Range("E41").Value = "40CEJICEJI"
